# "Remove Facility from entire app"
# The workbook's single data sheet has a "Facility" header in column C
# (with sample data "Facility name 1" in row 2). Delete that column
# entirely so everything to its right shifts left by one, and update
# the selection to match the post-delete cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Facility" column (column C) completely - this shifts all
# subsequent columns (D..AM) one place to the left and drops the
# now-unused "Facility" / "Facility name 1" shared strings.
$ws.Columns("C").Delete()

# Reset the view/selection to the top-left area of the sheet.
$ws.Range("C1").Select()
